$d = $word.ActiveDocument

# Locate the paragraph that begins with the Q26 question text
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Q26.*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the Q26 paragraph to replace"
}

# Remove the paragraph's contents (keep the paragraph mark / following empty para intact)
$d.Paragraphs($targetIndex).Range.Delete()

# Insert the replacement paragraphs in batches (the host has a limit on how many
# paragraphs a single InsertXML call can reliably insert at a point). Each batch is
# inserted at the paragraph-index anchor left behind by the previous batch, since the
# Range returned/updated by InsertXML only tracks the first inserted paragraph.

# Batch 1: Q41 - Q50
$xmlBatch0 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Q41</w:t></w:r><w:r><w:tab/><w:t>Subnetting is an important concept in network management. Can you explain what subnetting is and how it helps in organizing and optimizing network resources?</w:t></w:r></w:p><w:p><w:r><w:t>Q42</w:t></w:r><w:r><w:tab/><w:t>In the context of computer networks, what is the purpose of the Domain Name System (DNS)? How does DNS resolve domain names to IP addresses?</w:t></w:r></w:p><w:p><w:r><w:t>Q43</w:t></w:r><w:r><w:tab/><w:t>DHCP (Dynamic Host Configuration Protocol) simplifies the assignment of IP addresses. Can you explain how DHCP works and its advantages in network administration?</w:t></w:r></w:p><w:p><w:r><w:t>Q44</w:t></w:r><w:r><w:tab/><w:t>The TCP/IP protocol suite forms the basis of modern internet communication. Can you discuss the key protocols within this suite and their specific functions?</w:t></w:r></w:p><w:p><w:r><w:t>Q45</w:t></w:r><w:r><w:tab/><w:t>Network security is essential in today''s interconnected world. Can you explain the role of a firewall in a computer network and how it enhances security?</w:t></w:r></w:p><w:p><w:r><w:t>Q46</w:t></w:r><w:r><w:tab/><w:t>Network Address Translation (NAT) is commonly used to allow multiple devices to share a single public IP address. Can you explain how NAT works and its benefits?</w:t></w:r></w:p><w:p><w:r><w:t>Q47</w:t></w:r><w:r><w:tab/><w:t>Proxy servers act as intermediaries between clients and servers. Can you discuss the purpose of a proxy server in a network and how it affects communication?</w:t></w:r></w:p><w:p><w:r><w:t>Q48</w:t></w:r><w:r><w:tab/><w:t>Gateways enable communication between different networks. Could you explain the role of a gateway and the functions it performs in network connectivity?</w:t></w:r></w:p><w:p><w:r><w:t>Q49</w:t></w:r><w:r><w:tab/><w:t>Virtual Local Area Networks (VLANs) are used to segment networks. Can you explain how VLANs work and their advantages in network administration?</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>Q50</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:tab/><w:t>Can you discuss the differences between physical addresses (MAC addresses) and logical addresses (IP addresses)? How do they function in network communication?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$batchRange0 = $d.Paragraphs($targetIndex).Range
$batchRange0.InsertXML($xmlBatch0)
$targetIndex = $targetIndex + 10

# Batch 2: Q51 - Q98
$xmlBatch1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Q51</w:t></w:r><w:r><w:tab/><w:t>Network protocols govern the rules and formats for communication. Can you provide examples of commonly used network protocols and their respective purposes?</w:t></w:r></w:p><w:p><w:r><w:t>Q90</w:t></w:r><w:r><w:tab/><w:t>Describe the concept of network latency and its impact.</w:t></w:r></w:p><w:p><w:r><w:t>Q91</w:t></w:r><w:r><w:tab/><w:t>What is a routing table and how is it used?</w:t></w:r></w:p><w:p><w:r><w:t>Q92</w:t></w:r><w:r><w:tab/><w:t>Explain the principles of wireless network communication.</w:t></w:r></w:p><w:p><w:r><w:t>Q93</w:t></w:r><w:r><w:tab/><w:t>What is a LAN?</w:t></w:r></w:p><w:p><w:r><w:t>Q94</w:t></w:r><w:r><w:tab/><w:t>What is a WAN?</w:t></w:r></w:p><w:p><w:r><w:t>Q95</w:t></w:r><w:r><w:tab/><w:t>Define a hub.</w:t></w:r></w:p><w:p><w:r><w:t>Q96</w:t></w:r><w:r><w:tab/><w:t>Explain a switch.</w:t></w:r></w:p><w:p><w:r><w:t>Q97</w:t></w:r><w:r><w:tab/><w:t>What is a MAC?</w:t></w:r></w:p><w:p><w:r><w:t>Q98</w:t></w:r><w:r><w:tab/><w:t>What is an IP?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$batchRange1 = $d.Paragraphs($targetIndex).Range
$batchRange1.InsertXML($xmlBatch1)
$targetIndex = $targetIndex + 10

# Batch 3: Q99 - Q108
$xmlBatch2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Q99</w:t></w:r><w:r><w:tab/><w:t>Define DHCP briefly.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Q100</w:t></w:r><w:r><w:tab/><w:t>Explain TCP/IP protocol suite.</w:t></w:r></w:p><w:p><w:r><w:t>Q101</w:t></w:r><w:r><w:tab/><w:t>What is a socket?</w:t></w:r></w:p><w:p><w:r><w:t>Q102</w:t></w:r><w:r><w:tab/><w:t>Define SSL/TLS briefly.</w:t></w:r></w:p><w:p><w:r><w:t>Q103</w:t></w:r><w:r><w:tab/><w:t>What is a firewall?</w:t></w:r></w:p><w:p><w:r><w:t>Q104</w:t></w:r><w:r><w:tab/><w:t>Explain NAT briefly.</w:t></w:r></w:p><w:p><w:r><w:t>Q105</w:t></w:r><w:r><w:tab/><w:t>What is a VLAN?</w:t></w:r></w:p><w:p><w:r><w:t>Q106</w:t></w:r><w:r><w:tab/><w:t>What is latency?</w:t></w:r></w:p><w:p><w:r><w:t>Q107</w:t></w:r><w:r><w:tab/><w:t>Define ARP briefly.</w:t></w:r></w:p><w:p><w:r><w:t>Q108</w:t></w:r><w:r><w:tab/><w:t>Explain packet switching.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$batchRange2 = $d.Paragraphs($targetIndex).Range
$batchRange2.InsertXML($xmlBatch2)
$targetIndex = $targetIndex + 10

# Batch 4: Q109 - Q118
$xmlBatch3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Q109</w:t></w:r><w:r><w:tab/><w:t>What is a VPN?</w:t></w:r></w:p><w:p><w:r><w:t>Q110</w:t></w:r><w:r><w:tab/><w:t>Define bandwidth briefly.</w:t></w:r></w:p><w:p><w:r><w:t>Q111</w:t></w:r><w:r><w:tab/><w:t>How does a router function in a network infrastructure, and what role does it play in directing data packets across different networks by examining IP addresses and making routing decisions?</w:t></w:r></w:p><w:p><w:r><w:t>Q112</w:t></w:r><w:r><w:tab/><w:t>Explain the concept of subnetting in IP addressing and how it enables the division of a large network into smaller subnetworks, enhancing network efficiency and scalability.</w:t></w:r></w:p><w:p><w:r><w:t>Q113</w:t></w:r><w:r><w:tab/><w:t>Discuss the advantages and disadvantages of various network topologies, such as star, mesh, ring, and bus, in terms of their ability to handle traffic, fault tolerance, and ease of administration.</w:t></w:r></w:p><w:p><w:r><w:t>Q114</w:t></w:r><w:r><w:tab/><w:t>What are the primary differences between TCP (Transmission Control Protocol) and UDP (User Datagram Protocol) in terms of their reliability, connection-oriented vs. connectionless communication, and suitability for different applications?</w:t></w:r></w:p><w:p><w:r><w:t>Q115</w:t></w:r><w:r><w:tab/><w:t>How does a network firewall function as a barrier between an internal network and external networks, filtering network traffic based on predefined security rules to protect against unauthorized access and potential security threats?</w:t></w:r></w:p><w:p><w:r><w:t>Q116</w:t></w:r><w:r><w:tab/><w:t>Explain the role of DNS (Domain Name System) in the internet''s infrastructure, including how it translates domain names into IP addresses and facilitates the mapping of human-readable domain names to machine-readable IP addresses.</w:t></w:r></w:p><w:p><w:r><w:t>Q117</w:t></w:r><w:r><w:tab/><w:t>What is NAT (Network Address Translation), and how does it enable the usage of private IP addresses within a network while allowing communication with external networks by translating private IP addresses to public IP addresses?</w:t></w:r></w:p><w:p><w:r><w:t>Q118</w:t></w:r><w:r><w:tab/><w:t>Discuss the differences between IPv4 and IPv6 addressing schemes, including the differences in address format, the number of available addresses, and the features introduced in IPv6 to address the limitations of IPv4.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$batchRange3 = $d.Paragraphs($targetIndex).Range
$batchRange3.InsertXML($xmlBatch3)
$targetIndex = $targetIndex + 10
